$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
